# Fix property_category values that were incorrectly left as "land".
#
# Sheet "建物" (Building): column I (property_category) rows 2-8 -> "building"
$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 8; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

# Sheet "汽車" (Car): column H (property_category) rows 2-3 -> "car"
$wsCar = $wb.Worksheets.Item("汽車")
for ($r = 2; $r -le 3; $r++) {
    $wsCar.Cells.Item($r, 8).Value = "car"
}
